# Weekly GitHub Actions refresh of the cryptos list: update each coin's
# Price (D) and Volume(1h) (E) text values. These columns store pre-
# formatted strings (e.g. "36.271.83", "1.00", "  -1.48%  "), not real
# numbers, so force the cell format to Text before assigning, otherwise
# Excel would silently coerce number-looking strings (like "54.48") into
# numeric values and drop things like trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.271.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.041.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.00%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.48"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.37"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.01%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.906"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.31"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.340.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.84%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.046.76"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.48"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.210.68"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.95%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.91"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.19"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.80%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.11%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.90"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.88"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.77%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.63%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0598"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.53%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0900"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.04"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.21"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.45%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.82%  "

# Rows 44 and 45 swap places (Cronos <-> Aave) with refreshed price/volume.
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.62"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.94%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0904"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.401.53"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.55"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.95"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.14%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.87"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.10%  "
